$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column N ---
# This shifts the old N/O/P columns (Late / heading / Outstanding) one slot to
# the right, becoming O/P/Q, and leaves a new, empty column N in their place.
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns.Item(14).Insert()

# The newly inserted column picks up the width of the column immediately to
# its left (column M / "In Advance").
$wsSchedule.Columns.Item(14).ColumnWidth = $wsSchedule.Columns.Item(13).ColumnWidth

# --- Active sheet / selection bookkeeping ---
# Make "Repayment schedule" the active sheet/tab (previously "Transactions" was
# active) and move its selection to S8. "Transactions" keeps its own prior
# selection untouched, it simply stops being the active tab.
$wsSchedule.Activate()
$wsSchedule.Range("S8").Select()
